$wb = $excel.ActiveWorkbook

# --- Summary sheet: Unmet Demand Penalty + derived outputs ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 0.01
$wsSummary.Range("B6").Value = -263724.5473518896
$wsSummary.Range("B7").Value = 10477853.63860845
$wsSummary.Range("B8").Value = 27790152.75161
$wsSummary.Range("B10").Value = 1684280.80911358

# --- DG Dispatch: gains dispatch values that were previously unmet demand ---
$wsDG = $wb.Worksheets.Item("DG Dispatch")
$dgValues = @{
    "K2" = 220.0898510449805
    "L2" = 235.7664149699872
    "N2" = 229.4130635965909
    "P2" = 231.2329957552695
    "K3" = 137.841438974359
    "N3" = 131.3417120833333
    "O3" = 142.5962444444444
    "L4" = 134.8846762812383
    "M4" = 138.9257839476051
    "K5" = 220.0898510449805
    "N5" = 229.4130635965909
    "O5" = 230.0982114216867
    "P5" = 231.2329957552695
    "K6" = 137.841438974359
    "M6" = 142.1340339220183
    "P6" = 133.9744074143302
    "L7" = 134.8846762812383
    "M7" = 138.9257839476051
    "N7" = 127.6855444652332
    "O7" = 138.4565384518428
    "K8" = 220.0898510449805
    "L8" = 235.7664149699872
    "M8" = 230.3462332272727
    "O8" = 230.0982114216867
    "K9" = 137.841438974359
    "M9" = 142.1340339220183
    "O9" = 142.5962444444444
    "Q9" = 139.9817740860215
    "L10" = 134.8846762812383
    "N10" = 127.6855444652332
    "M11" = 230.3462332272727
    "N11" = 229.4130635965909
    "P11" = 231.2329957552695
    "K12" = 137.841438974359
    "L12" = 138.5543797798742
    "M12" = 142.1340339220183
    "N12" = 131.3417120833333
    "O12" = 142.5962444444444
    "P12" = 133.9744074143302
    "N13" = 127.6855444652332
    "O13" = 138.4565384518428
    "M14" = 230.3462332272727
    "N14" = 229.4130635965909
    "O14" = 230.0982114216867
    "P14" = 231.2329957552695
    "K15" = 137.841438974359
    "M15" = 142.1340339220183
    "N15" = 131.3417120833333
    "P15" = 133.9744074143302
    "Q15" = 139.9817740860215
    "M16" = 138.9257839476051
    "K17" = 220.0898510449805
    "L17" = 235.7664149699872
    "N17" = 229.4130635965909
    "P17" = 231.2329957552695
    "M18" = 142.1340339220183
    "P18" = 133.9744074143302
    "M20" = 230.3462332272727
    "L21" = 138.5543797798742
    "O21" = 142.5962444444444
    "P21" = 133.9744074143302
    "N22" = 127.6855444652332
    "L23" = 235.7664149699872
    "M23" = 230.3462332272727
    "P24" = 133.9744074143302
    "N25" = 127.6855444652332
    "O25" = 138.4565384518428
    "L26" = 235.7664149699872
    "M26" = 230.3462332272727
    "N26" = 229.4130635965909
    "N27" = 131.3417120833333
    "Q27" = 139.9817740860215
    "L28" = 134.8846762812383
    "K29" = 220.0898510449805
    "L29" = 235.7664149699872
    "M29" = 230.3462332272727
    "N29" = 229.4130635965909
    "P29" = 231.2329957552695
    "K30" = 137.841438974359
    "L30" = 138.5543797798742
    "M30" = 142.1340339220183
    "P30" = 133.9744074143302
    "L31" = 134.8846762812383
    "M31" = 138.9257839476051
    "N31" = 127.6855444652332
    "N32" = 229.4130635965909
    "K33" = 137.841438974359
    "L33" = 138.5543797798742
    "M33" = 142.1340339220183
    "P33" = 133.9744074143302
    "L34" = 134.8846762812383
    "M34" = 138.9257839476051
    "N34" = 127.6855444652332
    "O34" = 138.4565384518428
    "K35" = 220.0898510449805
    "N35" = 229.4130635965909
    "K36" = 137.841438974359
    "L36" = 138.5543797798742
    "M36" = 142.1340339220183
    "N36" = 131.3417120833333
    "P36" = 133.9744074143302
    "L37" = 134.8846762812383
    "M37" = 138.9257839476051
    "O37" = 138.4565384518428
    "K38" = 220.0898510449805
    "L38" = 235.7664149699872
    "N38" = 229.4130635965909
    "K39" = 137.841438974359
    "Q39" = 139.9817740860215
    "L40" = 134.8846762812383
    "M40" = 138.9257839476051
    "N40" = 127.6855444652332
    "M41" = 230.3462332272727
    "O41" = 230.0982114216867
    "P41" = 231.2329957552695
    "K42" = 137.841438974359
    "L42" = 138.5543797798742
    "Q42" = 139.9817740860215
    "L43" = 134.8846762812383
    "N43" = 127.6855444652332
    "O43" = 138.4565384518428
    "N44" = 229.4130635965909
    "L45" = 138.5543797798742
    "M45" = 142.1340339220183
    "Q45" = 139.9817740860215
    "M46" = 138.9257839476051
}
foreach ($ref in $dgValues.Keys) {
    $wsDG.Range($ref).Value = $dgValues[$ref]
}

# --- Unmet Demand: loses the same values (now served by DG) ---
$wsUD = $wb.Worksheets.Item("Unmet Demand")
$udRefs = @(
    "K2", "L2", "N2", "P2", "K3", "N3", "O3", "L4", "M4", "K5", "N5", "O5", "P5", "K6", "M6", "P6", "L7", "M7", "N7", "O7", "K8", "L8", "M8", "O8", "K9", "M9", "O9", "Q9", "L10", "N10", "M11", "N11", "P11", "K12", "L12", "M12", "N12", "O12", "P12", "N13", "O13", "M14", "N14", "O14", "P14", "K15", "M15", "N15", "P15", "Q15", "M16", "K17", "L17", "N17", "P17", "M18", "P18", "M20", "L21", "O21", "P21", "N22", "L23", "M23", "P24", "N25", "O25", "L26", "M26", "N26", "N27", "Q27", "L28", "K29", "L29", "M29", "N29", "P29", "K30", "L30", "M30", "P30", "L31", "M31", "N31", "N32", "K33", "L33", "M33", "P33", "L34", "M34", "N34", "O34", "K35", "N35", "K36", "L36", "M36", "N36", "P36", "L37", "M37", "O37", "K38", "L38", "N38", "K39", "Q39", "L40", "M40", "N40", "M41", "O41", "P41", "K42", "L42", "Q42", "L43", "N43", "O43", "N44", "L45", "M45", "Q45", "M46"
)
foreach ($ref in $udRefs) {
    $wsUD.Range($ref).Value = 0
}

# --- Household Surplus ---
$wsHS = $wb.Worksheets.Item("Household Surplus")
$hsValues = @{
    "B2" = 172147.4564623187
    "B3" = 182847.5994019398
    "B4" = 182987.3945782901
    "B5" = 171165.9949734709
    "B6" = 187632.7892677333
    "B7" = 150272.5348460527
    "B8" = 69227.33336802496
    "B9" = 87277.62269806072
    "B10" = 124436.7205980396
    "B11" = 220723.4332581452
    "B12" = 96997.22070799567
    "B13" = 132539.0888193136
    "B14" = 134834.408395709
    "B15" = 146147.4231337193
    "B16" = 69991.64848464866
}
foreach ($ref in $hsValues.Keys) {
    $wsHS.Range($ref).Value = $hsValues[$ref]
}

# --- Costs and Revenues: Total Operation Variable Costs (row4) & Total Profits (row6) ---
$wsCR = $wb.Worksheets.Item("Costs and Revenues")
$crValues = @{
    "B4" = 2424.612062849559
    "C4" = 2575.318301435772
    "D4" = 2577.287247581552
    "E4" = 2410.788661598182
    "F4" = 2642.715341799061
    "G4" = 2116.514575296517
    "H4" = 975.0328643383805
    "I4" = 1229.262291521983
    "J4" = 1752.629867578022
    "K4" = 3108.780750114724
    "L4" = 1366.158038140784
    "M4" = 1866.747729849489
    "N4" = 1899.076174587451
    "O4" = 2058.414410334074
    "P4" = 985.7978659809678
    "B6" = -54153.64424660709
    "C6" = -54153.64424660708
    "D6" = -54153.64424660709
    "E6" = -20526.04424660708
    "F6" = -20526.04424660708
    "G6" = -20526.04424660709
    "H6" = -20526.04424660708
    "I6" = -20526.04424660709
    "J6" = -20526.04424660708
    "K6" = -20526.04424660708
    "L6" = -20526.04424660708
    "M6" = -20526.04424660708
    "N6" = -20526.04424660708
    "O6" = -20526.04424660708
    "P6" = -20526.04424660708
}
foreach ($ref in $crValues.Keys) {
    $wsCR.Range($ref).Value = $crValues[$ref]
}

